# edit.ps1
# Adds a new data column ("28. 9. 2021") to both worksheets of the
# ZBP_10b_imunizace workbook, updates the trailing "aktualizace" date in the
# title rows, and lets Excel extend each sheet's used-range/dimension
# naturally as a side effect of writing the new cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data" -- percentages, new column AI, rows 1-75 (row 76 is the
# title-only footer row).
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

# Header: clone AH1's formatting (border/bold/alignment) onto AI1, then set
# the new header label.
$wsData.Range("AH1").Copy($wsData.Range("AI1"))
$wsData.Range("AI1").Value = "28. 9. 2021"

# Data rows 2-75: plain numeric values (no special style), exactly as the
# sibling AH column cells.
$aiValues = @{
    2 = 0.66
    3 = 0.05
    4 = 0.58
    5 = 0.05
    6 = 0.6
    7 = 0.07000000000000001
    8 = 0.77
    9 = 0.03
    10 = 0.61
    11 = 0.06
    12 = 0.65
    13 = 0.04
    14 = 0.79
    15 = 0.04
    16 = 0.53
    17 = 0.07000000000000001
    18 = 0.5600000000000001
    19 = 0.05
    20 = 0.75
    21 = 0.04
    22 = 0.7
    23 = 0.04
    24 = 0.86
    25 = 0.02
    26 = 0.88
    27 = 0.03
    28 = 0.43
    29 = 0.07000000000000001
    30 = 0.65
    31 = 0.04
    32 = 0.71
    33 = 0.06
    34 = 0.76
    35 = 0.03
    36 = 0.72
    37 = 0.05
    38 = 0.61
    39 = 0.05
    40 = 0.68
    41 = 0.03
    42 = 0.66
    43 = 0.06
    44 = 0.8
    45 = 0.04
    46 = 0.46
    47 = 0.06
    48 = 0.54
    49 = 0.07000000000000001
    50 = 0.74
    51 = 0.02
    52 = 0.63
    53 = 0.06
    54 = 0.66
    55 = 0.04
    56 = 0.64
    57 = 0.03
    58 = 0.8
    59 = 0.02
    60 = 0.7
    61 = 0.03
    62 = 0.77
    63 = 0.02
    64 = 0.67
    65 = 0.03
    66 = 0.79
    67 = 0.03
    68 = 0.59
    69 = 0.06
    70 = 0.7
    71 = 0.05
    72 = 0.7
    73 = 0.05
    74 = 0.64
    75 = 0.05
}
foreach ($row in $aiValues.Keys) {
    $wsData.Cells.Item([int]$row, 35).Value = $aiValues[$row]
}

# Footer title row: bump the "aktualizace" date.
$wsData.Range("A76").Value = "Život během pandemie, Imunizace, % respondentů celkově a ve skupinách, aktualizace 6. 10. 2021"

# ---------------------------------------------------------------------
# Sheet "pocetR" -- sample sizes, new column AH, rows 1-38 (row 39 is the
# title/footer row).
# ---------------------------------------------------------------------
$wsPocet = $wb.Worksheets.Item("pocetR")

# Header: clone AG1's formatting onto AH1, then set the new header label.
$wsPocet.Range("AG1").Copy($wsPocet.Range("AH1"))
$wsPocet.Range("AH1").Value = "28. 9. 2021"

# Data rows 2-38: plain numeric values.
$ahValues = @{
    2 = 1855
    3 = 456
    4 = 678
    5 = 721
    6 = 859
    7 = 631
    8 = 365
    9 = 441
    10 = 444
    11 = 250
    12 = 418
    13 = 187
    14 = 115
    15 = 180
    16 = 679
    17 = 613
    18 = 255
    19 = 897
    20 = 958
    21 = 237
    22 = 343
    23 = 317
    24 = 219
    25 = 336
    26 = 403
    27 = 965
    28 = 429
    29 = 216
    30 = 245
    31 = 238
    32 = 218
    33 = 265
    34 = 278
    35 = 584
    36 = 317
    37 = 325
    38 = 1213
}
foreach ($row in $ahValues.Keys) {
    $wsPocet.Cells.Item([int]$row, 34).Value = $ahValues[$row]
}

# Row 39 already carries an empty (but styled/used) cell in every column
# through AG; extend that same "touched" footer cell into the new AH
# column so the sheet's used range covers AH39 as well.
$wsPocet.Range("AH39").Formula = '=""'

# Footer title row: bump the "aktualizace" date.
$wsPocet.Range("A39").Value = "Život během pandemie, Imunizace, velikost dotázaného souboru celkově a ve skupinách, aktualizace 6. 10. 2021"
